# "Generate Report for Handoff" — refresh the localization-status report
# with the results of a new handoff cycle for da84d9ba-122b-4ba2-b0ec-f366cb6f244f.md.
#
# Overview!G6  "Latest HO Xliff Generate Date" -> 2016-09-03 02:45:01
# zh-cn!H6     "Latest Handoff Datetime"       -> 2016-09-03 02:44:56
# de-de!H6     "Latest Handoff Datetime"       -> 2016-09-03 02:45:01

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-09-03 02:45:01"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-09-03 02:44:56"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H6").Value = "2016-09-03 02:45:01"
